$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC01")
$ws.Activate()

$ws.Range("B1").Value = "Login endpoint"

$ws.Range("B1:F1").Select()
